# Update to user data
# - Notes sheet: update "Source: 40" -> "Source: 38"
# - Data sheet: add a new 2016 data row for each county (after its 2014 row)
# - Data-wide-value sheet: add a new "2016" column with the matching values

$wb = $excel.ActiveWorkbook

# --- 1. Notes sheet: fix the Source note ---
$wsNotes = $wb.Worksheets.Item("Notes")
$wsNotes.Range("A4").Value = "Source: 38"

# --- 2. Data sheet: insert one 2016 row per county, right after its 2014 row ---
$wsData = $wb.Worksheets.Item("Data")

# district_id|district_name|value_2016, in the same county order as the existing rows
$countyData2016 = @(
    "d18974|Baringo|152000",
    "d18975|Bomet|179000",
    "d18987|Bungoma|321000",
    "d18988|Busia|177000",
    "d18976|Elgeyo-Marakwet|99000",
    "d18955|Embu|164000",
    "d18965|Garissa|78000",
    "d18968|Homa Bay|224000",
    "d18956|Isiolo|34000",
    "d18991|Kajiado|250000",
    "d18989|Kakamega|392000",
    "d18977|Kericho|211000",
    "d18943|Kiambu|600000",
    "d18949|Kilifi|326000",
    "d18944|Kirinyaga|198000",
    "d18969|Kisii|291000",
    "d18970|Kisumu|284000",
    "d18957|Kitui|236000",
    "d18950|Kwale|174000",
    "d18978|Laikipia|135000",
    "d18951|Lamu|30000",
    "d18958|Machakos|328000",
    "d18959|Makueni|233000",
    "d18966|Mandera|111000",
    "d18960|Marsabit|62000",
    "d18961|Meru|393000",
    "d18971|Migori|233000",
    "d18952|Mombasa|397000",
    "d18946|Murang'a|323000",
    "d18964|Nairobi|1503000",
    "d18979|Nakuru|578000",
    "d18980|Nandi|202000",
    "d18981|Narok|223000",
    "d18962|Nithi|107000",
    "d18972|Nyamira|179000",
    "d18947|Nyandarua|191000",
    "d18948|Nyeri|271000",
    "d18982|Samburu|61000",
    "d18973|Siaya|246000",
    "d18953|Taita Taveta|102000",
    "d18954|Tana River|56000",
    "d18983|Trans-Nzoia|210000",
    "d18984|Turkana|246000",
    "d18985|Uasin Gishu|270000",
    "d18990|Vihiga|144000",
    "d18967|Wajir|69000",
    "d18986|West Pokot|119000"
)

$rowIdx = 4
foreach ($line in $countyData2016) {
    $parts = $line.Split("|")
    $wsData.Rows.Item($rowIdx).Insert()
    $wsData.Cells.Item($rowIdx, 1).Value = $parts[0]
    $wsData.Cells.Item($rowIdx, 2).Value = $parts[1]
    $wsData.Cells.Item($rowIdx, 3).Value = 2016
    $wsData.Cells.Item($rowIdx, 4).Value = [double]$parts[2]
    $rowIdx += 3
}

# --- 3. Data-wide-value sheet: add the "2016" column (D) ---
$wsWide = $wb.Worksheets.Item("Data-wide-value")
# Force the header to be stored as text (not a number) like the other year columns
$wsWide.Range("D1").NumberFormat = "@"
$wsWide.Range("D1").Value = "2016"
$wsWide.Range("D1").Style = "Normal"

$rowIdx = 2
foreach ($line in $countyData2016) {
    $parts = $line.Split("|")
    $wsWide.Cells.Item($rowIdx, 4).Value = [double]$parts[2]
    $rowIdx += 1
}
